$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "Untreated" group row is inserted right after the header row; the
# previously-existing group rows (Control/leisure/parallel combinations) all
# shift down by one row but otherwise keep their original labels and values.
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with the "Untreated" group's data.
$ws.Cells.Item(2, 1).Value = "Untreated"
$ws.Cells.Item(2, 2).Value = 8.2
$ws.Cells.Item(2, 3).Value = 28.85
$ws.Cells.Item(2, 4).Value = 31.3
$ws.Cells.Item(2, 5).Value = 32.29
